$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a brand-new row above current row 400, pushing rows 400..499 down to 401..500.
$ws.Rows.Item(400).Insert()

# Populate the newly inserted row 400 with a full record (mirrors the static
# columns of the old row 400, with fresh values for the changed columns).
$ws.Cells.Item(400, 1).Value = 3
$ws.Cells.Item(400, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(400, 3).Value = "Coquimbo"
$ws.Cells.Item(400, 4).Value = 44932
$ws.Cells.Item(400, 5).Value = 5
$ws.Cells.Item(400, 6).Value = 100112031
$ws.Cells.Item(400, 7).Value = "Poroto verde"
$ws.Cells.Item(400, 8).Value = "Magnum"
$ws.Cells.Item(400, 9).Value = "Primera"
$ws.Cells.Item(400, 10).Value = 73
$ws.Cells.Item(400, 11).Value = 28000
$ws.Cells.Item(400, 12).Value = 29000
$ws.Cells.Item(400, 13).Value = 28479
$ws.Cells.Item(400, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(400, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(400, 16).Value = 1139
$ws.Cells.Item(400, 17).Value = 25
$ws.Cells.Item(400, 18).Value = "Hortaliza"
